$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CasesTab query text in B2: remove the trailing Cohort column
$newQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
		WHERE s.clinical_study_designation IN ['UBC01'] and demo.sex in ['Female'] and demo.neutered_indicator in ['No']  OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@
$ws.Range("B2").Value = $newQuery

# Row 2 shrinks (autofit) now that the query text has one less line
$ws.Rows.Item(2).RowHeight = 290

# Selection moves from C2:F2 to B2 (and viewport scrolls back to A2)
[void]$ws.Range("B2").Select()
